$wb = $excel.ActiveWorkbook
$wb | Get-Member
